$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '38.660.74'
$ws.Cells.Item(2, 5).Value = '  +2.43%  '

$ws.Cells.Item(3, 4).Value = '2.097.19'
$ws.Cells.Item(3, 5).Value = '  +3.21%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).Formula = "'229.12"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.58%  '

$ws.Cells.Item(6, 4).Formula = "'0.615"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.62%  '

$ws.Cells.Item(7, 4).Formula = "'61.60"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +2.68%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 4).Formula = "'0.383"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +1.80%  '

$ws.Cells.Item(10, 4).Formula = "'0.0846"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +2.37%  '

$ws.Cells.Item(11, 4).Formula = "'0.104"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.46%  '

$ws.Cells.Item(12, 4).Value = '2.409.74'
$ws.Cells.Item(12, 5).Value = '  +3.33%  '

$ws.Cells.Item(13, 4).Formula = "'14.83"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +1.37%  '

$ws.Cells.Item(14, 4).Formula = "'22.42"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +6.47%  '

$ws.Cells.Item(15, 4).Formula = "'0.785"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +2.14%  '

$ws.Cells.Item(16, 4).Formula = "'5.47"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +5.42%  '

$ws.Cells.Item(17, 4).Value = '2.103.17'
$ws.Cells.Item(17, 5).Value = '  +3.40%  '

$ws.Cells.Item(18, 4).Value = '38.500.84'
$ws.Cells.Item(18, 5).Value = '  +2.09%  '

$ws.Cells.Item(19, 4).Formula = "'71.15"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +2.25%  '

$ws.Cells.Item(20, 5).Value = '  +2.68%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0835'
$ws.Cells.Item(21, 5).Value = '  +1.28%  '

$ws.Cells.Item(22, 4).Formula = "'226.41"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.17%  '

$ws.Cells.Item(23, 5).Value = '  -0.06%  '

$ws.Cells.Item(24, 5).Value = '  -0.66%  '

$ws.Cells.Item(25, 5).Value = '  +1.92%  '

$ws.Cells.Item(26, 4).Formula = "'170.49"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +2.06%  '

$ws.Cells.Item(27, 4).Formula = "'9.46"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.83%  '

$ws.Cells.Item(28, 4).Formula = "'0.133"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +2.27%  '

$ws.Cells.Item(29, 4).Formula = "'19.14"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +1.95%  '

$ws.Cells.Item(30, 5).Value = '  +7.47%  '

$ws.Cells.Item(31, 5).Value = '  -0.17%  '

$ws.Cells.Item(32, 5).Value = '  +2.88%  '

$ws.Cells.Item(33, 5).Value = '  +6.14%  '

$ws.Cells.Item(34, 4).Formula = "'4.49"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +2.53%  '

$ws.Cells.Item(35, 5).Value = '  +0.19%  '

$ws.Cells.Item(36, 4).Formula = "'6.58"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.87%  '

$ws.Cells.Item(37, 5).Value = '  +3.49%  '

$ws.Cells.Item(38, 4).Formula = "'3.56"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +4.20%  '

$ws.Cells.Item(39, 4).Formula = "'1.00"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.09%  '

$ws.Cells.Item(40, 4).Formula = "'18.76"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +3.34%  '

$ws.Cells.Item(41, 4).Value = '1.546.59'
$ws.Cells.Item(41, 5).Value = '  +0.72%  '

$ws.Cells.Item(42, 4).Formula = "'100.14"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +4.18%  '

$ws.Cells.Item(43, 4).Formula = "'0.0220"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +2.36%  '

$ws.Cells.Item(44, 4).Formula = "'2.83"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.35%  '

$ws.Cells.Item(45, 4).Formula = "'0.0915"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.72%  '

$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).Formula = "'7.66"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +8.02%  '

$ws.Cells.Item(47, 2).Value = 'FTXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(47, 4).Formula = "'4.17"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +3.42%  '

$ws.Cells.Item(48, 5).Value = '  +1.12%  '

$ws.Cells.Item(49, 5).Value = '  +3.30%  '

$ws.Cells.Item(50, 5).Value = '  +0.64%  '

$ws.Cells.Item(51, 4).Value = '2.294.43'
$ws.Cells.Item(51, 5).Value = '  +3.32%  '
